$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each target cell to Text format before writing so numeric-looking
# strings (e.g. "603.70") are preserved verbatim as text rather than being
# coerced into Excel numbers, then restore the default "Normal" style so no
# stray number-format style is left attached to the cell.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '70.381.26'
$c.Style = 'Normal'

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +0.92%  '
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.519.73'
$c.Style = 'Normal'

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -0.50%  '
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.18%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '603.70'
$c.Style = 'Normal'

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -0.87%  '
$c.Style = 'Normal'

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '174.90'
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +1.35%  '
$c.Style = 'Normal'

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.611'
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.80%  '
$c.Style = 'Normal'

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '3.512.56'
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -0.68%  '
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.23%  '
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.194'
$c.Style = 'Normal'

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -1.56%  '
$c.Style = 'Normal'

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '7.24'
$c.Style = 'Normal'

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +7.56%  '
$c.Style = 'Normal'

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.582'
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.32%  '
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '46.26'
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -3.21%  '
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000276'
$c.Style = 'Normal'

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -1.69%  '
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '4.080.63'
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -0.18%  '
$c.Style = 'Normal'

$c = $ws.Range('B16')
$c.NumberFormat = '@'
$c.Value = 'Polkadot'
$c.Style = 'Normal'

$c = $ws.Range('C16')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '8.32'
$c.Style = 'Normal'

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -1.22%  '
$c.Style = 'Normal'

$c = $ws.Range('B17')
$c.NumberFormat = '@'
$c.Value = 'BitcoinCash'
$c.Style = 'Normal'

$c = $ws.Range('C17')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '613.46'
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -1.01%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '3.513.12'
$c.Style = 'Normal'

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '70.373.80'
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +1.09%  '
$c.Style = 'Normal'

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +0.40%  '
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '17.28'
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -0.41%  '
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.876'
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -1.22%  '
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.07'
$c.Style = 'Normal'

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -19.64%  '
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '15.63'
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -1.70%  '
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '97.34'
$c.Style = 'Normal'

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.51%  '
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '3.72'
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -4.63%  '
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -0.03%  '
$c.Style = 'Normal'

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -2.97%  '
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '34.24'
$c.Style = 'Normal'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +2.31%  '
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.02'
$c.Style = 'Normal'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -3.51%  '
$c.Style = 'Normal'

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '8.16'
$c.Style = 'Normal'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -4.50%  '
$c.Style = 'Normal'

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.98'
$c.Style = 'Normal'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -5.46%  '
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '645.35'
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +12.95%  '
$c.Style = 'Normal'

$c = $ws.Range('B34')
$c.NumberFormat = '@'
$c.Value = 'NEARProtocol'
$c.Style = 'Normal'

$c = $ws.Range('C34')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '6.91'
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -1.27%  '
$c.Style = 'Normal'

$c = $ws.Range('B35')
$c.NumberFormat = '@'
$c.Value = 'Mantle'
$c.Style = 'Normal'

$c = $ws.Range('C35')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c.Style = 'Normal'

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.28'
$c.Style = 'Normal'

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -4.88%  '
$c.Style = 'Normal'

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.60'
$c.Style = 'Normal'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.Style = 'Normal'

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0997'
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -2.58%  '
$c.Style = 'Normal'

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '10.72'
$c.Style = 'Normal'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -1.06%  '
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0475'
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +6.21%  '
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '56.67'
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -1.42%  '
$c.Style = 'Normal'

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +0.18%  '
$c.Style = 'Normal'

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +1.29%  '
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.0₃0750'
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +5.32%  '
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '3.364.01'
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -0.84%  '
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.311'
$c.Style = 'Normal'

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -5.34%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.92'
$c.Style = 'Normal'

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +1.19%  '
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '32.25'
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -2.86%  '
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.56'
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -2.39%  '
$c.Style = 'Normal'

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.130'
$c.Style = 'Normal'

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.Style = 'Normal'

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +0.25%  '
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'

